$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 435 (existing rows 435:457 shift down to 438:460)
$ws.Rows("435:437").Insert()

# Static column values shared by every data row in this block
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"
$origen      = "Región de O'Higgins"

# New rows data: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, PrecioKg, KgUnidad
$newRows = @(
    @{ Row=435; D=44578; K="Candy White"; L="Extra (doble especial)"; M=230; N=16000; O=16000; P=16000; Q="$/bandeja 18 kilos granel"; S=889; T=18 },
    @{ Row=436; D=44578; K="Sun Rise";    L="Primera";                M=300; N=12000; O=12000; P=12000; Q="$/bandeja 18 kilos granel"; S=667; T=18 },
    @{ Row=437; D=44578; K="Venus";       L="Extra (doble especial)"; M=200; N=15000; O=15000; P=15000; Q="$/bandeja 15 kilos granel"; S=1000; T=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
